$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values (dependent formula cells recalc automatically)
$ws.Range("C4").Value = 70
$ws.Range("C5").Value = 0.02
$ws.Range("C15").Value = 2450

# Update the active selection to match the saved view state
$ws.Range("B5").Select()
